$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data values in column B (rows 1-6)
$ws.Range("B1").Value = 893
$ws.Range("B2").Value = 306
$ws.Range("B3").Value = 510
$ws.Range("B4").Value = 572
$ws.Range("B5").Value = 902
$ws.Range("B6").Value = 826

# Add new row of data (row 7)
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 710

# Update the chart: extend the series data ranges to include row 7
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$chart.SeriesCollection(1).Formula = "=SERIES(,'Sheet1'!`$A`$1:`$A`$7,'Sheet1'!`$B`$1:`$B`$7,1)"

# Clear the chart title text
$chart.ChartTitle.Text = ""
